$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "branch_name" header (column A) is removed; every other column
# (product_name, quantity, units, expire_date, price, description, status)
# shifts one slot to the left.
$ws.Range("A1").EntireColumn.Delete()

# After the deletion, the sheet's active selection moved to A2.
$ws.Range("A2").Select()

# Resize the workbook window (cosmetic, matches the saved view state).
$excel.ActiveWindow.Width = 19190
$excel.ActiveWindow.Height = 7520
